# Auto-committed on 2022/08/26 週五 17:20:49.47
#
# Adds two new field-definition rows (SEQ 13 & 14) to the "DBD" sheet of
# JcicZ045.xlsx:
#   Row 22: 13 | ActualFilingDate | 實際報送日期 | Decimald  | 8
#   Row 23: 14 | ActualFilingMark | 實際報送記號 | VARCHAR2  | 3
# and leaves the selection on C26 (matching the author's last cursor spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- Row 22 : SEQ 13, ActualFilingDate / 實際報送日期 / Decimald / 8 ---
$ws.Cells.Item(22, 1).Value = 13
$ws.Cells.Item(22, 2).Value = "ActualFilingDate"
$ws.Cells.Item(22, 3).Value = "實際報送日期"
$ws.Cells.Item(22, 4).Value = "Decimald"
$ws.Cells.Item(22, 5).Value = 8

# --- Row 23 : SEQ 14, ActualFilingMark / 實際報送記號 / VARCHAR2 / 3 ---
$ws.Cells.Item(23, 1).Value = 14
$ws.Cells.Item(23, 2).Value = "ActualFilingMark"
$ws.Cells.Item(23, 3).Value = "實際報送記號"
$ws.Cells.Item(23, 4).Value = "VARCHAR2"
$ws.Cells.Item(23, 5).Value = 3

# Move the active selection to C26, matching the saved view state.
$ws.Activate()
$ws.Range("C26").Select()
